$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (borders, style) of the last existing data row down
# onto the two new rows before filling in their values.
$ws.Range("A155:G155").Copy()
$ws.Range("A156:G157").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 156 - Thar (Najran region, south of the kingdom)
$ws.Range("A156").Value = "Thar"
$ws.Range("B156").Value = "Thar"
$ws.Range("C156").Value = "ثار"
$ws.Range("D156").Value = 17.981574999999999
$ws.Range("E156").Value = 44.105584999999998
$ws.Range("F156").Value = "منطقة نجران"
$ws.Range("G156").Value = "جنوب المملكة"

# Row 157 - Al Farshah (Asir region, south of the kingdom)
$ws.Range("A157").Value = "Al Farshah"
$ws.Range("B157").Value = "Al Farshah"
$ws.Range("C157").Value = "الفرشة"
$ws.Range("D157").Value = 17.753822
$ws.Range("E157").Value = 43.154803999999999
$ws.Range("F157").Value = "منطقة عسير"
$ws.Range("G157").Value = "جنوب المملكة"

# Update the selection to cover the full, now-larger table range
$ws.Range("A1:G157").Select()
